$wb = $excel.ActiveWorkbook

# --- Sheet "Sheet": append new work-study log entries (rows 161-176) ---
$ws = $wb.Worksheets.Item("Sheet")

$entries = @(
    "20 10:30>>> 0FD89CD000   David Schachner",
    "20 14:58>>> 0FD89CD000   David Schachner",
    "23 08:55>>> 0FD8AE8B60   Natalie Primus",
    "23 09:10>>> 0FD89CD000   David Schachner",
    "23 11:03>>> CFD88517C0   Joseph Giordano",
    "24 12:58>>> CFD893A460   Gary Tsai",
    "24 14:26>>> CFD893A460   Gary Tsai",
    "24 14:28>>> 8FD8A5E5E0   Anthony Avevor",
    "24 14:37>>> 8FD8A5E5E0   Anthony Avevor",
    "24 14:39>>> 8FD8A5E5E0   Anthony Avevor",
    "24 14:39>>> CFD893A460   Gary Tsai",
    "24 14:48>>> CFD8A948E0   Vanessa Cordero",
    "24 14:56>>> CFD893A460   Gary Tsai",
    "24 14:56>>> CFD893A460   Gary Tsai",
    "24 14:57>>> CFD893A460   Gary Tsai",
    "24 15:00>>> CFD893A460   Gary Tsai"
)

$startRow = 161
for ($i = 0; $i -lt $entries.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $entries[$i]
}
